$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify an existing entry's name (Mario Ordoñez Ramirez -> Mario Ordoñez)
$ws.Range("B130").Value = "Mario Ordoñez"

# Grab the table object so the new rows become part of Table1 (auto-expanding
# the table range, autofilter, dimension, etc.)
$tbl = $ws.ListObjects.Item(1)

$row133 = $tbl.ListRows.Add()
$row133.Range.Item(1, 1).Value = "Félix Romero"
$row133.Range.Item(1, 2).Value = "Kevin Jimenez"

$row134 = $tbl.ListRows.Add()
$row134.Range.Item(1, 1).Value = "Mario Ordoñez"
$row134.Range.Item(1, 2).Value = "Glenda Correa"

# Match formatting used elsewhere in the table: normal banded-row style for
# the new interior row, and the special bottom-border style for the new
# last row of the table.
$ws.Range("A2:B2").Copy() | Out-Null
$ws.Range("A133:B133").PasteSpecial(-4122) | Out-Null

$ws.Range("A132:B132").Copy() | Out-Null
$ws.Range("A134:B134").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Leave the cursor/view roughly where the real edit session left it.
$ws.Range("B135").Select() | Out-Null
